# Update column F (dSF) values for specific rows as per repull/recalculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 2
    10 = -1
    11 = 0
    12 = -3
    13 = -7
    14 = -2
    15 = -7
    16 = -2
    19 = -3
    20 = -2
    21 = 0
    24 = 3
    25 = 2
    26 = 3
    27 = 1
    33 = -1
    35 = 1
    37 = -9
    41 = -2
    43 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
